$wb = $excel.ActiveWorkbook

# 1. Bump version 1.8.1 -> 1.8.2 on the Metadata sheet (row "Version").
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.8.2"

# 2. Add a new "Include from ..." sheet for the new CodeSystem, mirroring
#    the existing "Include from identifierType" sheet layout/content.
#    Inserted after the last existing sheet, to keep it last in tab order.
$src = $wb.Worksheets.Item("Include from identifierType")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Include from Tipo Identificad"

$newSheet.Columns.Item(1).ColumnWidth = 30.703125
$newSheet.Columns.Item(2).ColumnWidth = 50.703125

$newSheet.Range("A1").Value = "Codes"
$newSheet.Range("A1").Style = $src.Range("A1").Style

$newSheet.Range("A2").Value = "All codes"
$newSheet.Range("A2").Style = $src.Range("A2").Style

$newSheet.Range("A3").Value = ""
$newSheet.Range("A3").Style = $src.Range("A3").Style
$newSheet.Range("B3").Value = ""
$newSheet.Range("B3").Style = $src.Range("B3").Style

$newSheet.Range("A4").Value = "System URI"
$newSheet.Range("A4").Style = $src.Range("A4").Style
$newSheet.Range("B4").Value = "https://hl7chile.cl/fhir/ig/clcore/CodeSystem/CSTipoIdentificador"
$newSheet.Range("B4").Style = $src.Range("B4").Style
